$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New experimental data rows (393-414) extracted from 10.1007/s11665-023-08899-x
# Alloy compositions (column B) - keep the exact spacing used by the source
# ---------------------------------------------------------------------------
$alloyFull   = "Al15Cr15Fe50Ni20"
$alloyB2sp   = "Al15Cr15Fe50Ni18B2 "
$alloyB4sp   = "Al15Cr15Fe50Ni16B4 "
$alloyB5sp   = "Al15Cr15Fe50Ni15B5 "
$alloyB6sp   = "Al15Cr15Fe50Ni14B6 "
$alloyB8     = "Al15Cr15Fe50Ni12B8"

$phaseFccBcc    = "FCC+BCC"
$phaseFccBccCr2B = "FCC+BCC+Cr2B"

# Column B (alloy name) - fill entire block first, column-major order matters
# for shared-string allocation.
$ws.Range("B393").Value = $alloyFull
$ws.Range("B394").Value = $alloyB2sp
$ws.Range("B395").Value = $alloyB4sp
$ws.Range("B396").Value = $alloyB5sp
$ws.Range("B397").Value = $alloyB6sp
$ws.Range("B398").Value = $alloyB8
$ws.Range("B399").Value = $alloyFull
$ws.Range("B400").Value = $alloyB2sp
$ws.Range("B401").Value = $alloyB4sp
$ws.Range("B402").Value = $alloyB5sp
$ws.Range("B403").Value = $alloyB6sp
$ws.Range("B404").Value = $alloyB8
$ws.Range("B405").Value = $alloyB4sp
$ws.Range("B406").Value = $alloyB5sp
$ws.Range("B407").Value = $alloyB6sp
$ws.Range("B408").Value = $alloyB8
$ws.Range("B409").Value = $alloyB4sp
$ws.Range("B410").Value = $alloyB5sp
$ws.Range("B411").Value = $alloyB6sp
$ws.Range("B412").Value = $alloyB8
$ws.Range("B413").Value = $alloyFull
$ws.Range("B414").Value = $alloyB2sp

# Column C (phase)
$ws.Range("C393").Value = $phaseFccBcc
$ws.Range("C394").Value = $phaseFccBcc
$ws.Range("C395").Value = $phaseFccBccCr2B
$ws.Range("C396").Value = $phaseFccBccCr2B
$ws.Range("C397").Value = $phaseFccBccCr2B
$ws.Range("C398").Value = $phaseFccBccCr2B
$ws.Range("C399").Value = $phaseFccBcc
$ws.Range("C400").Value = $phaseFccBcc
$ws.Range("C401").Value = $phaseFccBccCr2B
$ws.Range("C402").Value = $phaseFccBccCr2B
$ws.Range("C403").Value = $phaseFccBccCr2B
$ws.Range("C404").Value = $phaseFccBccCr2B
$ws.Range("C405").Value = $phaseFccBccCr2B
$ws.Range("C406").Value = $phaseFccBccCr2B
$ws.Range("C407").Value = $phaseFccBccCr2B
$ws.Range("C408").Value = $phaseFccBccCr2B
$ws.Range("C409").Value = $phaseFccBccCr2B
$ws.Range("C410").Value = $phaseFccBccCr2B
$ws.Range("C411").Value = $phaseFccBccCr2B
$ws.Range("C412").Value = $phaseFccBccCr2B
$ws.Range("C413").Value = $phaseFccBcc
$ws.Range("C414").Value = $phaseFccBcc

# Column D (processing route) - AAM for every new row
$ws.Range("D393:D414").Value = "AAM"

# Column F (property name)
$ws.Range("F393:F398").Value = "hardness"
$ws.Range("F399:F404").Value = "compressive yield stress"
$ws.Range("F405:F408").Value = "UCS"
$ws.Range("F409:F412").Value = "compressive ductility"
$ws.Range("F413:F414").Value = "minimum compressive ductility"

# Column G (method) - EXP for every new row
$ws.Range("G393:G414").Value = "EXP"

# Column I (temperature, K)
$ws.Range("I393:I414").Value = 298

# Column P (raw reported values backing the hardness/stress/UCS formulas)
$ws.Range("P393").Value = 433
$ws.Range("P394").Value = 485
$ws.Range("P395").Value = 492
$ws.Range("P396").Value = 575
$ws.Range("P397").Value = 537
$ws.Range("P398").Value = 498
$ws.Range("P399").Value = 1096
$ws.Range("P400").Value = 1243
$ws.Range("P401").Value = 1258
$ws.Range("P402").Value = 1330
$ws.Range("P403").Value = 1778
$ws.Range("P404").Value = 1912
$ws.Range("P405").Value = 2564
$ws.Range("P406").Value = 2595
$ws.Range("P407").Value = 2414
$ws.Range("P408").Value = 2070

# Column J - unit-converted values.
# Hardness block: Pa = P * 9807000 (kgf/mm^2 -> Pa), one shared formula group.
$ws.Range("J393:J398").Formula = "=P393*9807000"

# Compressive yield stress (399-404) + UCS (405-408): Pa = P * 1e6 (MPa -> Pa).
# Row 399 keeps its own (non-shared) formula, matching the source file, while
# 400-408 form a single shared-formula block.
$ws.Range("J399").Formula = "=P399*1000000"
$ws.Range("J400:J408").Formula = "=P400*1000000"

# Compressive ductility + minimum compressive ductility: plain % values.
$ws.Range("J409").Value = 42
$ws.Range("J410").Value = 31
$ws.Range("J411").Value = 23
$ws.Range("J412").Value = 16
$ws.Range("J413").Value = 80
$ws.Range("J414").Value = 80

# Column L (units)
$ws.Range("L393:L408").Value = "Pa"
$ws.Range("L409:L414").Value = "%"

# Column M (table reference id in the source paper)
$ws.Range("M393:M414").Value = "T3"

# Column N (DOI of the source paper)
$ws.Range("N393:N414").Value = "10.1007/s11665-023-08899-x"

# ---------------------------------------------------------------------------
# Drop the stray formatted-but-empty B/C/D placeholders at the very bottom of
# the sheet (rows 925-926), mirroring the author's final save.
# ---------------------------------------------------------------------------
$ws.Range("B925:D926").Clear()

# ---------------------------------------------------------------------------
# Leave the viewport parked where the author left it after entering the data.
# ---------------------------------------------------------------------------
$ws.Range("N418").Select()
